# This workbook is CRM test data: several cells hold digit-only phone
# numbers / dates / times that must remain *text* (shared-string) cells,
# not get auto-coerced into numbers or date serials the way a plain
# Range.Value assignment would. To force literal text while leaving the
# cell's existing style untouched, we compute the text via a formula in
# a scratch cell far outside the used range, copy it, and paste back
# only the value (xlPasteValues = -4163) onto the target cell.
$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$cellRef, [string]$text)
    $esc = $text.Replace('"', '""')
    $ws.Range("ZZ100").Formula = '="' + $esc + '"'
    $ws.Range("ZZ100").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("ZZ100").Clear()
}

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet1: refreshed record (new mobile numbers / run count / date /
# time / timestamp) plus a new "CurrentTime" column (header in AZ1, value
# in AZ2, replacing the old SMS-template message column). ---
Set-TextValue $ws1 "F2"  "4169897104"
Set-TextValue $ws1 "AE2" "2835679243"
Set-TextValue $ws1 "AT2" "9341419469"
Set-TextValue $ws1 "AX2" "5300619648"
Set-TextValue $ws1 "AK2" "3"
Set-TextValue $ws1 "N2"  "2024-03-06"
Set-TextValue $ws1 "O2"  "02:35:55 PM"
Set-TextValue $ws1 "P2"  "2024-03-06 07:16:02 PM"
Set-TextValue $ws1 "AZ1" "CurrentTime"
Set-TextValue $ws1 "AZ2" "CT. Wed, Mar 06, 2024 at 7:21 PM"
# AZ2 no longer needs the wrapped-text style the old long message used.
$ws1.Range("AZ2").WrapText = $false
Set-TextValue $ws1 "AC2" "2024-03-06"

# --- Sheets 2-4: same refreshed mobile numbers, no other columns touched. ---
Set-TextValue $ws2 "F2"  "4169897104"
Set-TextValue $ws2 "AE2" "2835679243"
Set-TextValue $ws2 "AT2" "9341419469"
Set-TextValue $ws2 "AX2" "5300619648"

Set-TextValue $ws3 "F2"  "4169897104"
Set-TextValue $ws3 "AE2" "2835679243"
Set-TextValue $ws3 "AT2" "9341419469"
Set-TextValue $ws3 "AX2" "5300619648"

Set-TextValue $ws4 "F2"  "4169897104"
Set-TextValue $ws4 "AE2" "2835679243"
Set-TextValue $ws4 "AT2" "9341419469"
Set-TextValue $ws4 "AX2" "5300619648"
